$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.814.84'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.646.04'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.500'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0628'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.18'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0842'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.49'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.648.04'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.17'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.58'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.820.61'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0738'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.00'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.58%  '
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.39'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +12.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.30'
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.38'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.01'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.11'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.71'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0513'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.32'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.00'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.292.21'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0174'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.540'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.828'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("E40").Value = '  +0.53%  '
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.36'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.797.37'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.62'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.21'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.48%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0518'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0980'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.11%  '
